$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.624.55"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.68%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.889.45"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.46%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "327.25"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("E6").Value = "  +0.09%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4606"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.36%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3875"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.70%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.81"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07887"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.84%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "21.79"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.67%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.894.13"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.096"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.726"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("E16").Value = "  -0.50%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "87.62"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("E18").Value = "  +0.09%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.00001006"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.25"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "28.644.86"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.347"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.05"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.132.89"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.057"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.06%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "154.88"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.40"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.18%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.889"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.24%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.964"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.26%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "118.63"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09359"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.41%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9268"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.311"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.342"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.79%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.264"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.65%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.05794"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.05%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "8.003"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.78%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.155"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.02073"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.27%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5702"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1798"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "9.803"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.5372"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.50%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "11.77"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.87%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.07151"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.41%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.176"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.19%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.846"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.118"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "112.85"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.479"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.52%  "
